$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Vega Monumental Concepción -
# Berenjena". It belongs chronologically right above the existing row 32
# (date 2021-09-13 / serial 44299), so insert a fresh row there and shift
# every row from 32..57 down to 33..58 (matches the diff exactly).
$ws.Rows.Item(32).Insert()

$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value = "Bíobío"
$ws.Cells.Item(32, 4).Value = 44554
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = 100112001
$ws.Cells.Item(32, 7).Value = "Berenjena"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 9000
$ws.Cells.Item(32, 12).Value = 10000
$ws.Cells.Item(32, 13).Value = 9500
$ws.Cells.Item(32, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(32, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value = 158
$ws.Cells.Item(32, 17).Value = 60
$ws.Cells.Item(32, 18).Value = "Hortaliza"
